$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list (coin price / 1h volume change %) values.
# Each target cell is forced to Text storage (NumberFormat "@") before the
# assignment so numeric-looking strings such as "305.61" or "-0.85%" are kept
# verbatim instead of being auto-converted to a number by the Value setter.
$updates = @{
    "D2" = "305.61"
    "E3" = "-0.85%"
    "D4" = "5.048"
    "E4" = "-0.95%"
    "D5" = "0.08040"
    "E5" = "-0.62%"
    "D6" = "1.907"
    "E6" = "-1.59%"
    "D7" = "4.154"
    "E7" = "-0.77%"
    "D8" = "7.780"
    "E8" = "0.31%"
    "E9" = "-0.70%"
    "D10" = "0.1280"
    "E10" = "-5.54%"
    "D11" = "0.1914"
    "E11" = "0.14%"
    "D12" = "0.09039"
    "E12" = "-1.40%"
    "D13" = "0.03450"
    "E13" = "1.23%"
    "D14" = "0.09860"
    "E14" = "0.31%"
    "D15" = "0.001404"
    "E15" = "-0.58%"
    "E16" = "7.47%"
    "D17" = "3.799"
    "E17" = "6.32%"
    "D18" = "3.395"
    "E18" = "13.83%"
    "D19" = "0.3420"
    "E19" = "-1.00%"
    "D20" = "0.1321"
    "E20" = "1.26%"
    "D21" = "5.185"
    "E21" = "5.73%"
    "D22" = "0.2391"
    "E22" = "-8.15%"
    "D23" = "0.04439"
    "E23" = "0.40%"
    "D24" = "0.001234"
    "E24" = "0.88%"
    "D25" = "0.004609"
    "E25" = "-4.26%"
    "E27" = "-3.89%"
    "E28" = "41.90%"
    "D39" = "0.01950"
    "E39" = "-2.86%"
    "D40" = "0.05427"
    "E40" = "10.44%"
    "D41" = "0.007615"
    "E41" = "-0.13%"
    "D42" = "0.01013"
    "E42" = "-1.05%"
    "D43" = "0.1353"
    "E43" = "-1.75%"
    "D44" = "0.002173"
    "E44" = "3.29%"
    "D45" = "0.009840"
    "E45" = "-10.67%"
    "D46" = "0.00006127"
    "E46" = "-4.23%"
    "E47" = "-0.03%"
    "D48" = "65.22"
    "E48" = "0.85%"
    "D49" = "0.001661"
    "E49" = "39.24%"
    "E50" = "-0.03%"
    "E51" = "-0.03%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
